$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.501.76"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").Value = "2.094.13"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  -3.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4428"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.58"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +15.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08941"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.153"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.26"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.12%  "

$ws.Range("D13").Value = "2.097.13"
$ws.Range("E13").Value = "  -1.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.691"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.686"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001122"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06611"
$ws.Range("D19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.253"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.80%  "

$ws.Range("D23").Value = "30.535.07"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("D26").Value = "2.343.09"
$ws.Range("E26").Value = "  -1.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.26"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.557"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.80"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.46"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1068"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.70%  "

$ws.Range("E33").Value = "  +6.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.155"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.903"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02556"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06816"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.479"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.59"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2257"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6869"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.251"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6322"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.196"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.626"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.237"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.243"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.97%  "

$ws.Range("E51").Value = "  -2.29%  "
